$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.663086891174316
$ws.Range("B1").Value = 2.225424289703369
$ws.Range("C1").Value = 3.240233421325684
$ws.Range("D1").Value = 4.499443531036377
$ws.Range("E1").Value = 0.6338819265365601
